# Team_Charter_Grid.xlsx — "testing and final clean up"
#
# Update the two pledged-contribution inputs on the "Costing" row (row 8)
# for the "Design Documents" task (columns D/E, merged) and the
# "Repository Schema" task (columns L/M, merged). The dependent SUM/
# difference formulas in row 19/20 recalculate automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D8").Value = 0.2
$ws.Range("L8").Value = 0.8

# Move the active selection to where the author left off editing.
$ws.Range("D18:E18").Select() | Out-Null
